# Update cryptos list values to match latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "299.73")
    # are not auto-coerced into Number cells, matching the source data
    # which always stores these as plain text.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '42.196.10'
Set-TextValue 'E2' '  -0.87%  '
Set-TextValue 'D3' '2.266.84'
Set-TextValue 'E3' '  -1.14%  '
Set-TextValue 'E4' '  -0.10%  '
Set-TextValue 'D5' '299.73'
Set-TextValue 'E5' '  -1.05%  '
Set-TextValue 'D6' '96.03'
Set-TextValue 'E6' '  -2.00%  '
Set-TextValue 'E7' '  -1.67%  '
Set-TextValue 'E8' '  -0.07%  '
Set-TextValue 'E9' '  -1.49%  '
Set-TextValue 'D10' '33.14'
Set-TextValue 'E10' '  -3.05%  '
Set-TextValue 'E11' '  +0.21%  '
Set-TextValue 'D12' '48.24'
Set-TextValue 'E12' '  -6.30%  '
Set-TextValue 'E13' '  +0.55%  '
Set-TextValue 'E14' '  -0.50%  '
Set-TextValue 'B15' 'Chainlink'
Set-TextValue 'C15' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D15' '15.54'
Set-TextValue 'E15' '  +0.02%  '
Set-TextValue 'B16' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C16' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D16' '2.618.75'
Set-TextValue 'E16' '  -1.30%  '
Set-TextValue 'D17' '2.294.13'
Set-TextValue 'E17' '  -1.28%  '
Set-TextValue 'D18' '0.783'
Set-TextValue 'E18' '  -2.62%  '
Set-TextValue 'D19' '42.114.29'
Set-TextValue 'E19' '  -1.03%  '
Set-TextValue 'D20' '11.68'
Set-TextValue 'E20' '  +2.12%  '
Set-TextValue 'D21' '0.0₃0888'
Set-TextValue 'E21' '  -0.83%  '
Set-TextValue 'D22' '5.98'
Set-TextValue 'E22' '  -0.77%  '
Set-TextValue 'D23' '66.34'
Set-TextValue 'E23' '  -3.13%  '
Set-TextValue 'D24' '234.94'
Set-TextValue 'E24' '  +0.43%  '
Set-TextValue 'D25' '1.97'
Set-TextValue 'E25' '  +0.20%  '
Set-TextValue 'D27' '2.45'
Set-TextValue 'E27' '  -2.18%  '
Set-TextValue 'E28' '  -4.07%  '
Set-TextValue 'D29' '168.52'
Set-TextValue 'E29' '  +3.60%  '
Set-TextValue 'D30' '2.06'
Set-TextValue 'E30' '  -4.81%  '
Set-TextValue 'D31' '9.17'
Set-TextValue 'E31' '  +0.44%  '
Set-TextValue 'E32' '  -2.48%  '
Set-TextValue 'E33' '  -0.22%  '
Set-TextValue 'D34' '4.89'
Set-TextValue 'E34' '  -2.07%  '
Set-TextValue 'D35' '4.53'
Set-TextValue 'E35' '  -1.47%  '
Set-TextValue 'D36' '16.59'
Set-TextValue 'E36' '  -1.73%  '
Set-TextValue 'E37' '  -4.34%  '
Set-TextValue 'E38' '  -3.61%  '
Set-TextValue 'E39' '  -3.15%  '
Set-TextValue 'D40' '0.0986'
Set-TextValue 'E40' '  -1.42%  '
Set-TextValue 'E41' '  -1.87%  '
Set-TextValue 'D42' '1.71'
Set-TextValue 'E42' '  -4.13%  '
Set-TextValue 'E43' '  -2.44%  '
Set-TextValue 'D44' '1.968.48'
Set-TextValue 'E44' '  -0.56%  '
Set-TextValue 'E45' '  -0.66%  '
Set-TextValue 'D46' '17.41'
Set-TextValue 'E46' '  -6.09%  '
Set-TextValue 'E47' '  -5.63%  '
Set-TextValue 'D48' '2.77'
Set-TextValue 'E48' '  -2.55%  '
Set-TextValue 'D49' '2.492.27'
Set-TextValue 'E49' '  -1.47%  '
Set-TextValue 'D50' '52.17'
Set-TextValue 'E50' '  -5.21%  '
Set-TextValue 'E51' '  -0.53%  '
